$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix grammar / expand ranges in the DESIGNATOR column (B)
$ws.Range("B12").Value = "R1, R4, R5, R7, R14,R15,R16,R17, R24, R26, R28,R32,R33,R34"
$ws.Range("B8").Value = "D6,D7,D8"
$ws.Range("B21").Value = "R18,R19,R20,R21,R25,R27,R29"

# Update the active selection to match where the last edit was made
$ws.Range("B21").Select()

$wb.Save()
